$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The app under test changed its pizza-vendor label and its result wording;
# re-enter the row values to reflect the new app output.
$ws.Range("A1").Value = "pizza"
$ws.Range("B1").Value = "Pizza Pizza"
$ws.Range("C1").Value = "PASSED"
$ws.Range("A2").Value = "chap"
$ws.Range("B2").Value = "Bikanervala"
$ws.Range("C2").Value = "PASSED"

$ws.Range("B6").Select()
